$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-118 shift down to 11-119.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the new record.
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44496
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112003
$ws.Range("G10").Value = "Ajo"
$ws.Range("H10").Value = "Chino"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 18000
$ws.Range("N10").Value = '$/caja 10 kilos'
$ws.Range("O10").Value = "China"
$ws.Range("P10").Value = 1800
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"
